# Horas_Equipo.xlsx edit: add "Run DataLake" row (tipo_proyecto) to general table,
# adjust the percentage split for Oscar Chero's "Run Datamart y cierres" row,
# and fix an inconsistent cell style on F8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proyectos")

# --- 1. Normalize F8's stray style (it carried an unused "applyFont" alignment
#        style that duplicates the plain center-aligned style used elsewhere).
#        Resetting Bold explicitly collapses it back onto the shared style.
$ws.Range("F8").Font.Bold = $false

# --- 2. Row 21 (Oscar Chero / Run Datamart y cierres): split changes 50% -> 40%
$ws.Range("F21").Formula = "=40*4*40%"
$ws.Range("G21:O21").Formula = "=40*4*40%"

# --- 2b. Row 22 (Alexander Atencio / Run ODS), before the insert below shifts
#         it to row 23: split changes 50% -> 60%, to account for the rebalance.
$ws.Range("F22").Formula = "=40*4*60%"
$ws.Range("G22:O22").Formula = "=40*4*60%"

# --- 3. Insert a brand-new row right after row 21 for the new "Run DataLake"
#        project line (everything below shifts down by one row, so the row we
#        just edited above is now row 23).
$ws.Rows.Item(22).Insert()

# Match formatting (center alignment) of the surrounding data rows for the
# newly inserted, still-blank row.
$ws.Range("A22:Q22").HorizontalAlignment = $ws.Range("A23:Q23").HorizontalAlignment
$ws.Range("A22:Q22").VerticalAlignment = $ws.Range("A23:Q23").VerticalAlignment

# --- 4. Populate the new row: Oscar Chero, Staff, Run, Run DataLake, 20%
$ws.Range("A22").Value = "Oscar Chero"
$ws.Range("B22").Value = "Staff"
$ws.Range("C22").Value = "Run"
$ws.Range("D22").Value = "Run DataLake"
$ws.Range("E22").Formula = "=SUM(F22:O22)"
$ws.Range("F22").Formula = "=40*4*20%"
$ws.Range("G22:O22").Formula = "=40*4*20%"

# --- 5. Restore the selection to where the editor left off.
$ws.Range("E22").Select()
